# TestData.xlsx update: "Atualização de projeto, ajustes."
#
# - Cadastro!B2        "Wilkerbn22"      -> "Wilkerbn504"
# - Produtos!C12       "P(11,0) V(11,1)" -> "C(11,0) V(11,1)"
# - Produtos!C13       "P(12,0) V(12,1)" -> "C(12,0) V(12,1)"
# - Produtos!C14       "P(13,0) V(13,1)" -> "C(13,0) V(13,1)"
# - Active sheet moves from "Produtos" (index 2) to "Cadastro" (index 1),
#   with the selection on Cadastro left at B2 and the stale selection on
#   Produtos moved to B18.

$wb = $excel.ActiveWorkbook

$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsProdutos = $wb.Worksheets.Item("Produtos")

# Update the shared-string backed cell values. Produtos' rows are updated
# before Cadastro's so the newly introduced shared strings land in the same
# relative order as the target workbook.
$wsProdutos.Range("C12").Value = "C(11,0) V(11,1)"
$wsProdutos.Range("C13").Value = "C(12,0) V(12,1)"
$wsProdutos.Range("C14").Value = "C(13,0) V(13,1)"

$wsCadastro.Range("B2").Value = "Wilkerbn504"

# Leave a stale selection behind on Produtos (it was the previously active
# sheet) and then switch the active sheet/selection to Cadastro, matching
# the bookViews activeTab="1" / per-sheet tabSelected + selection in the
# target workbook.
[void]$wsProdutos.Activate()
[void]$wsProdutos.Range("B18").Select()

[void]$wsCadastro.Activate()
[void]$wsCadastro.Range("B2").Select()
